$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 104.181816
$ws.Range("I33").Value = 96.22221999999999
$ws.Range("K33").Value = 96.22221999999999
$ws.Range("M33").Value = 132.77778
$ws.Range("H76").Value = 7722.5557
$ws.Range("J76").Value = 7928.5713
$ws.Range("L76").Value = 7928.5713
$ws.Range("N76").Value = -8558.5713
$ws.Range("H79").Value = 7722.5557
$ws.Range("J79").Value = 7928.5713
$ws.Range("L79").Value = 7928.5713
$ws.Range("N79").Value = -10112.5713
$ws.Range("H135").Value = 1586.909
$ws.Range("J135").Value = 1650
$ws.Range("L135").Value = 14850
$ws.Range("N135").Value = -19920
$ws.Range("H137").Value = 1999.5
$ws.Range("I137").Value = 1999
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 5997
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -3447
$ws.Range("N137").Value = -11100

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 140
$ws.Range("I5").Value = 140
$ws.Range("K5").Value = 140
$ws.Range("M5").Value = -28
$ws.Range("H63").Value = 5663.25
$ws.Range("I63").Value = 1326.5
$ws.Range("K63").Value = 1326.5
$ws.Range("M63").Value = -640.5
$ws.Range("H66").Value = 5663.25
$ws.Range("I66").Value = 1326.5
$ws.Range("K66").Value = 6632.5
$ws.Range("M66").Value = -3200.5
$ws.Range("H74").Value = 2489.1538
$ws.Range("I74").Value = 2085.9
$ws.Range("K74").Value = 2085.9
$ws.Range("M74").Value = -1211.9
$ws.Range("H77").Value = 2489.1538
$ws.Range("I77").Value = 2085.9
$ws.Range("K77").Value = 10429.5
$ws.Range("M77").Value = -6061.5
$ws.Range("H132").Value = 2530.7058
$ws.Range("I132").Value = 2386.6155
$ws.Range("K132").Value = 7159.8465
$ws.Range("M132").Value = -4629.8465

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 140
$ws.Range("I4").Value = 140
$ws.Range("K4").Value = 140
$ws.Range("M4").Value = -25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5100.643
$ws.Range("J31").Value = 9623.5
$ws.Range("L31").Value = 9623.5
$ws.Range("N31").Value = -10213.5
$ws.Range("H34").Value = 5100.643
$ws.Range("J34").Value = 9623.5
$ws.Range("L34").Value = 9623.5
$ws.Range("N34").Value = -10027.5
$ws.Range("H58").Value = 2854.2
$ws.Range("I58").Value = 2222.7273
$ws.Range("J58").Value = 3626
$ws.Range("K58").Value = 2222.7273
$ws.Range("L58").Value = 3626
$ws.Range("M58").Value = -2019.7273
$ws.Range("N58").Value = -4032
$ws.Range("H100").Value = 40000
$ws.Range("I100").Value = 40000
$ws.Range("K100").Value = 40000
$ws.Range("M100").Value = -38918
$ws.Range("H107").Value = 1007.3333
$ws.Range("I107").Value = 595.55
$ws.Range("K107").Value = 595.55
$ws.Range("M107").Value = 1324.45
$ws.Range("H133").Value = 29999
$ws.Range("J133").Value = 29999
$ws.Range("L133").Value = 29999
$ws.Range("N133").Value = -35059
$ws.Range("H136").Value = 2854.2
$ws.Range("I136").Value = 2222.7273
$ws.Range("J136").Value = 3626
$ws.Range("K136").Value = 6668.1819
$ws.Range("L136").Value = 10878
$ws.Range("M136").Value = -4118.1819
$ws.Range("N136").Value = -15978

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 111125620
$ws.Range("J64").Value = 166686670
$ws.Range("L64").Value = 500060010
$ws.Range("N64").Value = -500060550
$ws.Range("H67").Value = 111125620
$ws.Range("J67").Value = 166686670
$ws.Range("L67").Value = 500060010
$ws.Range("N67").Value = -500061882
$ws.Range("H129").Value = 1550
$ws.Range("J129").Value = 3000
$ws.Range("L129").Value = 9000
$ws.Range("N129").Value = -19000
$ws.Range("H131").Value = 2326.8447
$ws.Range("J131").Value = 2385.037
$ws.Range("L131").Value = 7155.110999999999
$ws.Range("N131").Value = -17235.111

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3091.5293
$ws.Range("J132").Value = 5753
$ws.Range("L132").Value = 17259
$ws.Range("N132").Value = -22319

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 403
$ws.Range("I12").Value = 403
$ws.Range("K12").Value = 403
$ws.Range("M12").Value = -233
$ws.Range("H16").Value = 10278.4
$ws.Range("I16").Value = 9898
$ws.Range("J16").Value = 10849
$ws.Range("K16").Value = 9898
$ws.Range("L16").Value = 10849
$ws.Range("M16").Value = -9728
$ws.Range("N16").Value = -11189
$ws.Range("H46").Value = 2011.4615
$ws.Range("I46").Value = 1540
$ws.Range("J46").Value = 2306.125
$ws.Range("K46").Value = 1540
$ws.Range("L46").Value = 2306.125
$ws.Range("M46").Value = -1352
$ws.Range("N46").Value = -2682.125
$ws.Range("H122").Value = 2614.3635
$ws.Range("I122").Value = 2189.875
$ws.Range("K122").Value = 6569.625
$ws.Range("M122").Value = -4119.625
$ws.Range("H125").Value = 287500
$ws.Range("J125").Value = 287500
$ws.Range("L125").Value = 287500
$ws.Range("N125").Value = -297340
$ws.Range("H132").Value = 211190.8
$ws.Range("J132").Value = 22500
$ws.Range("L132").Value = 67500
$ws.Range("N132").Value = -72560

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 50011
$ws.Range("J20").Value = 50011
$ws.Range("L20").Value = 50011
$ws.Range("N20").Value = -50491
$ws.Range("H28").Value = 24909.5
$ws.Range("J28").Value = 24909.5
$ws.Range("L28").Value = 24909.5
$ws.Range("N28").Value = -25605.5
$ws.Range("H30").Value = 22505
$ws.Range("I30").Value = 5000
$ws.Range("J30").Value = 40010
$ws.Range("K30").Value = 5000
$ws.Range("L30").Value = 40010
$ws.Range("M30").Value = -4893
$ws.Range("N30").Value = -40224
$ws.Range("H33").Value = 27843.715
$ws.Range("J33").Value = 27984.334
$ws.Range("L33").Value = 27984.334
$ws.Range("N33").Value = -28484.334
$ws.Range("H36").Value = 27843.715
$ws.Range("J36").Value = 27984.334
$ws.Range("L36").Value = 27984.334
$ws.Range("N36").Value = -28484.334
$ws.Range("H37").Value = 50029
$ws.Range("J37").Value = 50029
$ws.Range("L37").Value = 50029
$ws.Range("N37").Value = -50435
$ws.Range("H45").Value = 25434.285
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("H126").Value = 1718.2667
$ws.Range("I126").Value = 1289.5
$ws.Range("J126").Value = 3433.3333
$ws.Range("K126").Value = 3868.5
$ws.Range("L126").Value = 10299.9999
$ws.Range("M126").Value = -1398.5
$ws.Range("N126").Value = -15239.9999
$ws.Range("H132").Value = 2866.4
$ws.Range("I132").Value = 2695.1428
$ws.Range("J132").Value = 3266
$ws.Range("K132").Value = 8085.428400000001
$ws.Range("L132").Value = 9798
$ws.Range("M132").Value = -5555.428400000001
$ws.Range("N132").Value = -14858
$ws.Range("H136").Value = 4611.9165
$ws.Range("I136").Value = 3305.0908
$ws.Range("K136").Value = 9915.2724
$ws.Range("M136").Value = -7365.2724
$ws.Range("M45").ClearContents()
